$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New, alphabetically sorted list of "cargo" values for B2:B56
$cargos = @(
    "Agente Comercial De Nacionales",
    "Analista Contable",
    "Analista De Documentos Y Facturador",
    "Analista De Procesos",
    "Analista De Reportes Y Documentos",
    "Aprendiz Sena",
    "Asistente De Compras",
    "Asistente De Gerencia",
    "Asistente De Logistica Y Transportes",
    "Asistente De Publicidad, Mercadeo Y Comunicaciones",
    "Asistente De Soporte Tecnico Ti",
    "Asistente De Tesoreria",
    "Asistente Tecnico",
    "Auxiliar Administrativo",
    "Auxiliar Contable",
    "Auxiliar De Almacen",
    "Auxiliar De Bodega",
    "Auxiliar De Logistica Y Operaciones",
    "Auxiliar De Recibo",
    "Auxiliar De Recursos Humanos",
    "Auxiliar De Servicios Generales",
    "Auxiliar Logistico Y Conductor",
    "Comprador",
    "Conductor",
    "Contador",
    "Coordinador  De Exportaciones Y Documentos",
    "Coordinador Administrativo Y Financiero",
    "Coordinador Contable",
    "Coordinador De Compras",
    "Coordinador De Logistica Y Transportes",
    "Coordinador De Operaciones",
    "Coordinador De Pagos",
    "Coordinador De Recibo",
    "Coordinador De Recursos Humanos",
    "Coordinador De Sst Y Procesos",
    "Coordinador Jr Aeropuerto",
    "Coordinador Tecnico Y De Compras ",
    "Gerente Administrativo y financiero",
    "Gerente Comercial Publicidad Y Mercadeo",
    "Gerente De Operaciones ",
    "Gerente General",
    "Gerente Tecnico Y De Compras",
    "Jefe de Tesoreria",
    "Operario",
    "Pasante Sst",
    "Patinador",
    "Recepcionista",
    "Representante De Ventas Interno",
    "Subgerente",
    "Supernumerario De Operaciones ",
    "Supervisor De Calidad",
    "Supervisor De Linea",
    "Supervisor De Produccion",
    "Supervisor De Produccion",
    "Supervisor De Trazabilidad"
)

for ($i = 0; $i -lt $cargos.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $cargos[$i]
}

# Update selection to match the author's final view (B7 active cell)
$ws.Range("B7").Select()
